$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Formatting first (creates / reuses cellXfs entries without
#    touching the shared-string table).
# ------------------------------------------------------------------

# E1 header cell gets the same (bold, bordered, centered) style as D1
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)

# New E5:E10 cells get the plain bordered/centered style used elsewhere
$ws.Range("B5").Copy()
$ws.Range("E5:E10").PasteSpecial(-4122)

# New E2:E4 cells (will hold "Paid") get the same plain style too
$ws.Range("B2").Copy()
$ws.Range("E2:E4").PasteSpecial(-4122)

# Date column A2:A8 needs a short-date number format (built-in id 14)
# plus the existing border/centered alignment. Apply to A2 then copy
# the resulting format down so only one new style entry is created.
$ws.Range("A2").NumberFormat = "mm-dd-yy"
$ws.Range("A2").Copy()
$ws.Range("A3:A8").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# 2) Text values, entered in the order needed to reproduce the
#    original shared-string table order: Cash, Sir, NEFT, Paid.
# ------------------------------------------------------------------
$ws.Range("C2").Value = "Cash"
$ws.Range("D2").Value = "Sir"
$ws.Range("C5").Value = "NEFT"
$ws.Range("E2").Value = "Paid"
$ws.Range("E3").Value = "Paid"
$ws.Range("E4").Value = "Paid"

# ------------------------------------------------------------------
# 3) Numeric / date / formula values.
# ------------------------------------------------------------------
$ws.Range("A2").Value = 45183
$ws.Range("A3").Value = 45189
$ws.Range("A4").Value = 45190
$ws.Range("A5").Value = 45194
$ws.Range("A6").Value = 45197
$ws.Range("A7").Value = 45199
$ws.Range("A8").Value = 45201

$ws.Range("B2").Value = 1000
$ws.Range("B5").Formula = "=900*3-B2"

# ------------------------------------------------------------------
# 4) Sheet view / selection.
# ------------------------------------------------------------------
$ws.Range("D5").Select()

# ------------------------------------------------------------------
# 5) Page setup (portrait orientation).
# ------------------------------------------------------------------
$ws.PageSetup.Orientation = 1
